# "edit v parser.c = exp misto R_RIGID. added some more semantic rules"
#
# The grammar table on Sheet1 has one row per production rule, with columns:
#   A = left-hand-side nonterminal, B = "::=" (shared formula), C = right-hand side.
# Row 14 held the rule  R_RIGID ::= = EXP  which became redundant once the
# parser started using "= EXP" directly (e.g. the R_FLEX / ALL_AFTER_ID rules),
# so that row is removed and every row below it shifts up one place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "R_RIGID ::= = EXP" rule (row 14); Excel shifts rows
# 15..59 up to 14..58 and fixes up the shared "::=" formulas automatically.
$ws.Rows(14).Delete()

# Match the author's final selection/viewport on the sheet.
$ws.Range("D17").Select()
